$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 in the "Journal de travail" table was previously an empty
# (but formatted) row. Fill it in with a new entry, matching the
# formatting already used by the row above (row 24).

# 1) Copy the formatting (number formats / styles / wrap text, etc.)
#    from row 24 down onto row 25 so the new values pick up the same
#    cell styles (date format, time format, centered+wrapped text...).
$ws.Range("E24:M24").Copy()
$ws.Range("E25:M25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Fill in the actual data for the new journal entry.
$ws.Range("E25").Value = 44265
$ws.Range("F25").Value = 0.75
$ws.Range("G25").Value = 0.79166666666666663

# Re-enter the duration formula so it recalculates against the new
# Heure Début / Heure fin values instead of keeping the pasted (stale)
# cached result.
$ws.Range("H25").Formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),"""",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])"

$ws.Range("I25").Value = "Développement"
$ws.Range("J25").Value = "Correction de bug + rajout de fonctionalité"
$ws.Range("K25").Value = "Domicile"
$ws.Range("L25").Value = "Correction d'érreur d'affichage pour la grille ainsi que vérifier si un bateau ce trouve sur la case des coordonnées et changer l'affichage en correspondance"

# 3) The long description in L25 wraps onto many lines, so the row
#    needs to grow tall enough to show all of it.
$ws.Rows("25").RowHeight = 158.4

# 4) Update the view: the user scrolled down a couple of rows and
#    moved the active selection from J25 to K26.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("K26").Select()

Write-Output "Row 25 filled in and view updated."
